# Commit: "Included a comment about .NET WCF on slide 3"
#
# Slide 3's bullet list currently reads (first bullet):
#   " Tjenester eksponert over ren HTTP"
# The edit splits that bullet's single run so the sentence keeps going,
# contrasting plain HTTP with SOAP (.NET WCF):
#   " Tjenester eksponert over ren HTTP i motsetning til over SOAP (.NET WCF)"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# Find the "TextBox 4" shape (holds the bullet list) by name, rather than
# assuming a fixed shape index.
$sh = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    if ($s.Shapes.Item($i).Name -eq "TextBox 4") {
        $sh = $s.Shapes.Item($i)
        break
    }
}
if ($sh -eq $null) {
    $sh = $s.Shapes.Item(3)
}

$tr = $sh.TextFrame.TextRange
$para1 = $tr.Paragraphs(1, 1)

# Original paragraph text (35 chars): " Tjenester eksponert over ren HTTP"
#   chars 1-26  -> " Tjenester eksponert over "
#   chars 27-30 -> "ren "
#   chars 31-34 -> "HTTP"
# Re-set the "ren " substring so it becomes its own run …
$runRen = $para1.Characters(27, 4)
$runRen.Text = "ren "

# … and extend the trailing "HTTP" substring into its own run with the new
# comment appended after it.
$runHttp = $para1.Characters(31, 4)
$runHttp.Text = "HTTP i motsetning til over SOAP (.NET WCF)"
